# The underlying commit swaps the contents of ppt/theme/theme1.xml (used by
# the slide master -> the deck's real visual theme, originally the green
# "Integral" scheme) and ppt/theme/theme2.xml (used only by the notes master,
# originally the default "Office Theme" scheme).
#
# Both theme parts share an identical fontScheme and fmtScheme; the only
# meaningful content difference between "Integral" and "Office Theme" is the
# 12-colour clrScheme. We reproduce the swap by rewriting the slide master's
# theme colours (the part of the object model this host exposes) from the
# "Integral" palette to the "Office Theme" palette.

function ConvertTo-VbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$tcs = $m.Theme.ThemeColorScheme

# Target palette: the "Office Theme" colours that used to live in theme2.xml,
# now moved onto the slide master's theme (theme1.xml), in the
# dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink (Colors(1..12)) order.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = ConvertTo-VbaRgb $officeThemeColors[$i - 1]
}
